$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3666.6667
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 3666.6667
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 11000.0001
$ws.Range("N69").Value = -12748.0001
$ws.Range("M69").ClearContents()

$ws.Range("H72").Value = 3666.6667
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 3666.6667
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 33000.0003
$ws.Range("N72").Value = -41736.0003
$ws.Range("M72").ClearContents()

$ws.Range("H81").Value = 33108.332
$ws.Range("J81").Value = 33108.332
$ws.Range("L81").Value = 33108.332
$ws.Range("N81").Value = -35104.332

$ws.Range("H84").Value = 33108.332
$ws.Range("J84").Value = 33108.332
$ws.Range("L84").Value = 99324.99600000001
$ws.Range("N84").Value = -109308.996

$ws.Range("H92").Value = 551.06665
$ws.Range("I92").Value = 519
$ws.Range("J92").Value = 1000
$ws.Range("K92").Value = 519
$ws.Range("L92").Value = 1000
$ws.Range("M92").Value = 729
$ws.Range("N92").Value = -3496

$ws.Range("H111").Value = 1849.3
$ws.Range("I111").Value = 2025.8
$ws.Range("J111").Value = 1672.8
$ws.Range("K111").Value = 6077.4
$ws.Range("L111").Value = 5018.4
$ws.Range("M111").Value = -3010.4
$ws.Range("N111").Value = -11152.4

$ws.Range("H116").Value = 9093162
$ws.Range("I116").Value = 20001618
$ws.Range("J116").Value = 2782.8333
$ws.Range("K116").Value = 20001618
$ws.Range("L116").Value = 2782.8333
$ws.Range("M116").Value = -19998176
$ws.Range("N116").Value = -9666.8333

$ws.Range("H132").Value = 2606.5
$ws.Range("I132").Value = 2552.15
$ws.Range("K132").Value = 7656.450000000001
$ws.Range("M132").Value = -5126.450000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1260.6923
$ws.Range("I16").Value = 1048.875
$ws.Range("J16").Value = 1599.6
$ws.Range("K16").Value = 1048.875
$ws.Range("L16").Value = 1599.6
$ws.Range("M16").Value = -761.875
$ws.Range("N16").Value = -2173.6

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H92").Value = 48257
$ws.Range("J92").Value = 48257
$ws.Range("L92").Value = 48257
$ws.Range("N92").Value = -53249

$ws.Range("H98").Value = 26000
$ws.Range("J98").Value = 26000
$ws.Range("L98").Value = 26000
$ws.Range("N98").Value = -30492

$ws.Range("H113").Value = 1260.6923
$ws.Range("I113").Value = 1048.875
$ws.Range("J113").Value = 1599.6
$ws.Range("K113").Value = 1048.875
$ws.Range("L113").Value = 1599.6
$ws.Range("M113").Value = 1121.125
$ws.Range("N113").Value = -5939.6

$ws.Range("H115").Value = 30000
$ws.Range("J115").Value = 30000
$ws.Range("L115").Value = 30000
$ws.Range("N115").Value = -32350

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("N17").ClearContents()

$ws.Range("H34").Value = 488.25
$ws.Range("I34").Value = 285
$ws.Range("J34").Value = 556
$ws.Range("K34").Value = 855
$ws.Range("L34").Value = 1668
$ws.Range("M34").Value = -771
$ws.Range("N34").Value = -1836

$ws.Range("H39").Value = 1694.75
$ws.Range("I39").Value = 500
$ws.Range("J39").Value = 2889.5
$ws.Range("K39").Value = 1500
$ws.Range("L39").Value = 8668.5
$ws.Range("M39").Value = -1206
$ws.Range("N39").Value = -9256.5

$ws.Range("H55").Value = 3405.5833
$ws.Range("I55").Value = 1000
$ws.Range("J55").Value = 3886.7
$ws.Range("K55").Value = 3000
$ws.Range("L55").Value = 11660.1
$ws.Range("N55").Value = -12014.1
$ws.Range("M55").Value = -2823

$ws.Range("H107").Value = 342.91666
$ws.Range("J107").Value = 275.6
$ws.Range("L107").Value = 826.8000000000001
$ws.Range("N107").Value = -4666.8

$ws.Range("H124").Value = 12199.8
$ws.Range("I124").Value = 999
$ws.Range("J124").Value = 15000
$ws.Range("K124").Value = 2997
$ws.Range("L124").Value = 45000
$ws.Range("M124").Value = 1913
$ws.Range("N124").Value = -54820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 781.1429000000001
$ws.Range("I107").Value = 745.7222
$ws.Range("K107").Value = 745.7222
$ws.Range("M107").Value = 1174.2778

$ws.Range("H113").Value = 1694.3684
$ws.Range("I113").Value = 1833.3334
$ws.Range("J113").Value = 1668.3125
$ws.Range("K113").Value = 1833.3334
$ws.Range("L113").Value = 1668.3125
$ws.Range("M113").Value = 336.6666
$ws.Range("N113").Value = -6008.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4429.923
$ws.Range("I7").Value = 4486.125
$ws.Range("J7").Value = 4340
$ws.Range("K7").Value = 4486.125
$ws.Range("L7").Value = 4340
$ws.Range("M7").Value = -4374.125
$ws.Range("N7").Value = -4564

$ws.Range("H61").Value = 5166.6665
$ws.Range("I61").Value = 7000
$ws.Range("J61").Value = 1500
$ws.Range("K61").Value = 7000
$ws.Range("L61").Value = 1500
$ws.Range("M61").Value = -6798
$ws.Range("N61").Value = -1904

$ws.Range("H68").Value = 5400.6
$ws.Range("I68").Value = 5500
$ws.Range("J68").Value = 5334.3335
$ws.Range("K68").Value = 5500
$ws.Range("L68").Value = 5334.3335
$ws.Range("M68").Value = -4751
$ws.Range("N68").Value = -6832.3335

$ws.Range("H71").Value = 5400.6
$ws.Range("I71").Value = 5500
$ws.Range("J71").Value = 5334.3335
$ws.Range("K71").Value = 27500
$ws.Range("L71").Value = 26671.6675
$ws.Range("M71").Value = -23756
$ws.Range("N71").Value = -34159.6675

$ws.Range("H100").Value = 4859.533
$ws.Range("I100").Value = 9032.166999999999
$ws.Range("J100").Value = 2077.7778
$ws.Range("K100").Value = 9032.166999999999
$ws.Range("L100").Value = 2077.7778
$ws.Range("M100").Value = -8491.166999999999
$ws.Range("N100").Value = -3159.7778

$ws.Range("H113").Value = 5166.6665
$ws.Range("I113").Value = 7000
$ws.Range("J113").Value = 1500
$ws.Range("K113").Value = 7000
$ws.Range("L113").Value = 1500
$ws.Range("M113").Value = -4830
$ws.Range("N113").Value = -5840

$ws.Range("H126").Value = 4429.923
$ws.Range("I126").Value = 4486.125
$ws.Range("J126").Value = 4340
$ws.Range("K126").Value = 13458.375
$ws.Range("L126").Value = 13020
$ws.Range("M126").Value = -10988.375
$ws.Range("N126").Value = -17960

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4270

$ws.Range("H65").Value = 4270

$ws.Range("H81").Value = 63809.5
$ws.Range("I81").Value = 63809.5
$ws.Range("K81").Value = 127619
$ws.Range("M81").Value = -126558

$ws.Range("H84").Value = 63809.5
$ws.Range("I84").Value = 63809.5
$ws.Range("K84").Value = 638095
$ws.Range("M84").Value = -632791

$ws.Range("H113").Value = 468.30768
$ws.Range("I113").Value = 535.55554
$ws.Range("J113").Value = 317
$ws.Range("K113").Value = 1606.66662
$ws.Range("L113").Value = 951
$ws.Range("M113").Value = 563.33338
$ws.Range("N113").Value = -5291

$ws.Range("H117").Value = 36704.5
$ws.Range("J117").Value = 36704.5
$ws.Range("L117").Value = 36704.5
$ws.Range("N117").Value = -45882.5
